$wb = $excel.ActiveWorkbook
$wsA = $wb.Worksheets.Item(1)

# --- Create sheets B, C, D after A ---
$wsB = $wb.Worksheets.Add($null, $wsA)
$wsB.Name = "B"

$wsC = $wb.Worksheets.Add($null, $wsB)
$wsC.Name = "C"

$wsD = $wb.Worksheets.Add($null, $wsC)
$wsD.Name = "D"

# --- Sheet B data ---
$wsB.Range("A1").Value = "curl de quadriceps"
$wsB.Range("B1").Value = 0
$wsB.Range("C1").Value = 60

$wsB.Range("A2").Value = "silla de aductores"
$wsB.Range("B2").Value = 0
$wsB.Range("C2").Value = 30

$wsB.Range("A3").Value = "sentadillas"
$wsB.Range("B3").Value = 1
$wsB.Range("C3").Value = 110

$wsB.Range("A4").Value = "subidas laterales al cajon"
$wsB.Range("B4").Value = 0
$wsB.Range("C4").Value = 0

$wsB.Range("A5").Value = "sentadillas bulgaras"
$wsB.Range("B5").Value = 0
$wsB.Range("C5").Value = 0

$wsB.Range("A6").Value = "pantorillas"
$wsB.Range("B6").Value = 0
$wsB.Range("C6").Value = 10

# --- Sheet C data ---
$wsC.Range("A1").Value = "press de hombros"
$wsC.Range("B1").Value = 72.5
$wsC.Range("C1").Value = 1

$wsC.Range("A2").Value = "vuelo frontal"
$wsC.Range("B2").Value = 9
$wsC.Range("C2").Value = 0

$wsC.Range("A3").Value = "vuelo lateral"
$wsC.Range("B3").Value = 9
$wsC.Range("C3").Value = 0

$wsC.Range("A4").Value = "vuelo posterior"
$wsC.Range("B4").Value = 9
$wsC.Range("C4").Value = 0

$wsC.Range("A5").Value = "biceps inclinado"
$wsC.Range("B5").Value = 9
$wsC.Range("C5").Value = 0

$wsC.Range("A6").Value = "biceps martillo"
$wsC.Range("B6").Value = 9
$wsC.Range("C6").Value = 0

# --- Sheet D data ---
$wsD.Range("A1").Value = "curl de femorales"
$wsD.Range("B1").Value = 5
$wsD.Range("C1").Value = 0

$wsD.Range("A2").Value = "sillon de abeductores"
$wsD.Range("B2").Value = 10
$wsD.Range("C2").Value = 0

$wsD.Range("A3").Value = "peso muerto"
$wsD.Range("B3").Value = 130
$wsD.Range("C3").Value = 1

$wsD.Range("A4").Value = "patadas laterales"
$wsD.Range("B4").Value = 2
$wsD.Range("C4").Value = 0

$wsD.Range("A5").Value = "hip thrust"
$wsD.Range("B5").Value = 80
$wsD.Range("C5").Value = 0

$wsD.Range("A6").Value = "pantorillas"
$wsD.Range("B6").Value = 10
$wsD.Range("C6").Value = 0

# --- Selections ---
$wsA.Range("A7").Select() | Out-Null
$wsB.Range("C7").Select() | Out-Null
$wsC.Range("C7").Select() | Out-Null
$wsD.Range("C7").Select() | Out-Null

# --- Activate sheet D (makes it the tabSelected / active tab) ---
$wsD.Activate() | Out-Null
